# Loan RBI, Variable Instalments
# - Insert a new (blank) column N on the "Repayment schedule" sheet, which
#   shifts the former N/O/P columns ("Late", "heading", "Outstanding") one
#   column to the right (-> O/P/Q) and widens the grid from A1:P9 to A1:Q9.
# - The newly inserted column inherits the width of the column to its left
#   (M), matching Excel's native "Insert Column" behaviour.
# - Switch the active sheet/selection from "Transactions" to
#   "Repayment schedule", selecting cell K13 there.

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the old column N ("Late").
[void]$schedule.Columns("N").Insert()

# Match the width Excel would copy in from the column immediately to the
# left (M) when a column is inserted.
$schedule.Columns("N").ColumnWidth = $schedule.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet (was "Transactions") and move
# the selection to K13.
[void]$schedule.Activate()
[void]$schedule.Range("K13").Select()
